$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "65.642.59"
$ws.Range("E2").Value = "  +0.73%  "

Set-TextValue "D3" "3.395.01"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue "D5" "561.26"
$ws.Range("E5").Value = "  -0.26%  "

Set-TextValue "D6" "175.85"
$ws.Range("E6").Value = "  +0.49%  "

Set-TextValue "D7" "0.631"
$ws.Range("E7").Value = "  +0.85%  "

Set-TextValue "D8" "3.388.16"
$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  +5.15%  "

Set-TextValue "D11" "0.634"
$ws.Range("E11").Value = "  +0.32%  "

Set-TextValue "D12" "53.56"
$ws.Range("E12").Value = "  -2.30%  "

Set-TextValue "D13" "0.0000278"
$ws.Range("E13").Value = "  +0.13%  "

Set-TextValue "D14" "9.20"
$ws.Range("E14").Value = "  +0.54%  "

Set-TextValue "D15" "3.939.12"
$ws.Range("E15").Value = "  -0.11%  "

Set-TextValue "D16" "18.25"
$ws.Range("E16").Value = "  -0.66%  "

Set-TextValue "D17" "3.406.01"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("E18").Value = "  +0.25%  "

Set-TextValue "D19" "65.681.84"
$ws.Range("E19").Value = "  +0.84%  "

Set-TextValue "D20" "11.83"
$ws.Range("E20").Value = "  -0.85%  "

Set-TextValue "D21" "0.998"
$ws.Range("E21").Value = "  +0.21%  "

Set-TextValue "D22" "479.37"
$ws.Range("E22").Value = "  +1.28%  "

Set-TextValue "D23" "4.95"
$ws.Range("E23").Value = "  -0.36%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "89.81"
$ws.Range("E24").Value = "  +3.68%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D25" "14.30"
$ws.Range("E25").Value = "  +3.34%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D26" "4.10"
$ws.Range("E26").Value = "  -1.17%  "

Set-TextValue "D27" "2.90"
$ws.Range("E27").Value = "  +0.64%  "

Set-TextValue "D28" "10.62"
$ws.Range("E28").Value = "  -2.45%  "

Set-TextValue "D29" "8.71"
$ws.Range("E29").Value = "  -1.53%  "

Set-TextValue "D30" "31.30"
$ws.Range("E30").Value = "  +1.84%  "

Set-TextValue "D31" "6.57"
$ws.Range("E31").Value = "  -1.93%  "

Set-TextValue "D32" "63.61"
$ws.Range("E32").Value = "  +5.33%  "

Set-TextValue "D33" "11.44"
$ws.Range("E33").Value = "  -0.89%  "

Set-TextValue "D34" "572.44"
$ws.Range("E34").Value = "  -2.21%  "

$ws.Range("E35").Value = "  -0.83%  "

$ws.Range("E36").Value = "  -0.01%  "

Set-TextValue "D37" "3.67"
$ws.Range("E37").Value = "  +3.81%  "

$ws.Range("E38").Value = "  +0.43%  "

Set-TextValue "D39" "35.90"
$ws.Range("E39").Value = "  -0.31%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D40" "0.0₃0747"
$ws.Range("E40").Value = "  -0.62%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D41" "0.373"
$ws.Range("E41").Value = "  -0.35%  "

Set-TextValue "D42" "3.082.98"
$ws.Range("E42").Value = "  -1.04%  "

Set-TextValue "D43" "2.80"
$ws.Range("E43").Value = "  -2.45%  "

$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("E45").Value = "  +0.09%  "

Set-TextValue "D46" "2.44"
$ws.Range("E46").Value = "  -3.03%  "

$ws.Range("E47").Value = "  -1.46%  "

Set-TextValue "D48" "0.999"
$ws.Range("E48").Value = "  -0.05%  "

Set-TextValue "D49" "140.38"
$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("E50").Value = "  +0.21%  "

Set-TextValue "D51" "8.45"
$ws.Range("E51").Value = "  +1.11%  "
